# daily auto push: 2026-01-29 09:56 UTC
# A new reading was recorded for 2026/01/29 (Thursday) at hour 17 that had
# been missing from the log. Insert it in chronological order as row 735,
# which pushes the existing rows 735-776 down to 736-777 (matching how the
# sheet's existing entries are naturally sorted by date/time).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 735..776 down to 736..777, opening up row 735 for the new entry.
$ws.Rows(735).Insert()

# Column A ("日付") is stored as plain text (e.g. "2026/12/29"), not a real
# date. Assigning a date-shaped string straight to .Value makes Excel coerce
# it into a date serial, so build it as a text formula first and then paste
# the computed value back over itself - that collapses it to a literal
# shared-string cell with no lingering formula or special number format,
# matching every other row in the column.
$ws.Range("A735").Formula = "=""2026/01/29"""
$ws.Range("A735").Copy()
$ws.Range("A735").PasteSpecial(-4163)
$excel.CutCopyMode = $false

$ws.Range("B735").Value = "木"
$ws.Range("C735").Value = 17
$ws.Range("D735").Value = 201
